$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = '2026-01-29'
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = '기술'
$ws.Cells.Item(2, 3).Value = '"맥락 분석해 오디오 자동 생성"…NC AI, 바르코 사운드 출시'
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '2026-01-29'
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = 'https://n.news.naver.com/mnews/article/421/0008741237?sid=105'

$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = '2026-01-29'
$ws.Cells.Item(3, 1).ClearFormats()
$ws.Cells.Item(3, 2).Value = '기술'
$ws.Cells.Item(3, 3).Value = '[AI픽] NC AI, 사운드 생성 인공지능 ''바르코 사운드'' 출시'
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2026-01-29'
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = 'https://n.news.naver.com/mnews/article/001/0015872911?sid=105'

$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = '2026-01-29'
$ws.Cells.Item(4, 1).ClearFormats()
$ws.Cells.Item(4, 2).Value = '정책'
$ws.Cells.Item(4, 3).Value = '우리은행 "공공기관 AI 전환 지원"…관련 협회와 업무협약'
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '2026-01-29'
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = 'https://n.news.naver.com/mnews/article/001/0015873297?sid=101'

$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = '2026-01-29'
$ws.Cells.Item(5, 1).ClearFormats()
$ws.Cells.Item(5, 2).Value = '산업'
$ws.Cells.Item(5, 3).Value = '작년 AI 적용 등 혁신의료기기 45개 지정…전년대비 1.5배'
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '2026-01-29'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = 'https://n.news.naver.com/mnews/article/001/0015873367?sid=105'

$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = '2026-01-29'
$ws.Cells.Item(6, 1).ClearFormats()
$ws.Cells.Item(6, 2).Value = '기업'
$ws.Cells.Item(6, 3).Value = '오픈AI에 맞불 놓은 구글, 월 1만1000원 저가 AI 요금제 출시'
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '2026-01-28'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = 'https://n.news.naver.com/mnews/article/003/0013734299?sid=105'

$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = '2026-01-29'
$ws.Cells.Item(7, 1).ClearFormats()
$ws.Cells.Item(7, 2).Value = '기업'
$ws.Cells.Item(7, 3).Value = '충남 천안에 80㎿급 AI 데이터센터 건립…투자협약 체결'
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '2026-01-29'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = 'https://n.news.naver.com/mnews/article/001/0015873655?sid=105'

$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = '2026-01-29'
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(8, 2).Value = '산업'
$ws.Cells.Item(8, 3).Value = '"AI 활용 가장 활발했다"…작년 혁신의료기기 45개 지정'
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '2026-01-29'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = 'https://n.news.naver.com/mnews/article/003/0013736629?sid=102'

$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = '2026-01-29'
$ws.Cells.Item(9, 1).ClearFormats()
$ws.Cells.Item(9, 2).Value = '정책'
$ws.Cells.Item(9, 3).Value = '우리은행, 공공기관 AI 전환 지원 나선다'
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2026-01-29'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = 'https://n.news.naver.com/mnews/article/277/0005713728?sid=101'

$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = '2026-01-29'
$ws.Cells.Item(10, 1).ClearFormats()
$ws.Cells.Item(10, 2).Value = '정부(과기부)'
$ws.Cells.Item(10, 3).Value = '트릴리온랩스, 국내 첫 확산 기반 트랜스포머 모델 개발'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '2026-01-29'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = 'https://n.news.naver.com/mnews/article/030/0003394619?sid=105'

$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = '2026-01-29'
$ws.Cells.Item(11, 1).ClearFormats()
$ws.Cells.Item(11, 2).Value = '정부(과기부)'
$ws.Cells.Item(11, 3).Value = '국내 연구자 호라이즌 유럽 참여 본격화…7개 과제 수주'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '2026-01-29'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = 'https://n.news.naver.com/mnews/article/001/0015874276?sid=105'

$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = '2026-01-29'
$ws.Cells.Item(12, 1).ClearFormats()
$ws.Cells.Item(12, 2).Value = '정부(과기부)'
$ws.Cells.Item(12, 3).Value = '과기정통부·중기부, AI스타트업 성장전략 설명회 개최'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '2026-01-29'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = 'https://www.itbiznews.com/news/articleView.html?idxno=202823'

$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = '2026-01-29'
$ws.Cells.Item(13, 1).ClearFormats()
$ws.Cells.Item(13, 2).Value = '정부(과기부)'
$ws.Cells.Item(13, 3).Value = '과기정통부, 양자 분야 첫 마스터플랜 공개…양자인력 1만명·기업 2000개...'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '2026-01-29'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = 'https://www.dailysecu.com/news/articleView.html?idxno=204778'

$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = '2026-01-29'
$ws.Cells.Item(14, 1).ClearFormats()
$ws.Cells.Item(14, 2).Value = '정부(과기부)'
$ws.Cells.Item(14, 3).Value = '배경훈 부총리 "AI시대, 양자역할 중요...투자 확대 가속화"'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2026-01-29'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = 'https://n.news.naver.com/mnews/article/014/0005470355?sid=105'
